# Applies the data update described in the commit diff to Sheets/Valefor_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 1932
$ws.Range("I62").Value = 1800
$ws.Range("J62").Value = 1998
$ws.Range("K62").Value = 1800
$ws.Range("L62").Value = 1998
$ws.Range("M62").Value = -1176
$ws.Range("N62").Value = -3246

# Row 65
$ws.Range("H65").Value = 1932
$ws.Range("I65").Value = 1800
$ws.Range("J65").Value = 1998
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 9990
$ws.Range("M65").Value = -5880
$ws.Range("N65").Value = -16230

# Row 74
$ws.Range("H74").Value = 4269.522
$ws.Range("I74").Value = 5058.25
$ws.Range("J74").Value = 3409.0908
$ws.Range("K74").Value = 5058.25
$ws.Range("L74").Value = 3409.0908
$ws.Range("M74").Value = -4122.25
$ws.Range("N74").Value = -5281.0908

# Row 76
$ws.Range("H76").Value = 2727.0908
$ws.Range("I76").Value = 2742.5715
$ws.Range("J76").Value = 2700
$ws.Range("K76").Value = 2742.5715
$ws.Range("L76").Value = 2700
$ws.Range("M76").Value = -2427.5715
$ws.Range("N76").Value = -3330

# Row 77
$ws.Range("H77").Value = 4269.522
$ws.Range("I77").Value = 5058.25
$ws.Range("J77").Value = 3409.0908
$ws.Range("K77").Value = 25291.25
$ws.Range("L77").Value = 17045.454
$ws.Range("M77").Value = -20611.25
$ws.Range("N77").Value = -26405.454

# Row 79
$ws.Range("H79").Value = 2727.0908
$ws.Range("I79").Value = 2742.5715
$ws.Range("J79").Value = 2700
$ws.Range("K79").Value = 2742.5715
$ws.Range("L79").Value = 2700
$ws.Range("M79").Value = -1650.5715
$ws.Range("N79").Value = -4884

# Row 106
$ws.Range("H106").Value = 3883.5
$ws.Range("I106").Value = 3500
$ws.Range("J106").Value = 4075.25
$ws.Range("K106").Value = 3500
$ws.Range("L106").Value = 4075.25
$ws.Range("M106").Value = -2869
$ws.Range("N106").Value = -5337.25

# Row 125
$ws.Range("H125").Value = 4678.6665
$ws.Range("J125").Value = 6518
$ws.Range("L125").Value = 58662
$ws.Range("N125").Value = -63582

# Row 129
$ws.Range("H129").Value = 1207.6046
$ws.Range("I129").Value = 321.375
$ws.Range("J129").Value = 1410.1714
$ws.Range("K129").Value = 964.125
$ws.Range("L129").Value = 4230.5142
$ws.Range("M129").Value = 4035.875
$ws.Range("N129").Value = -14230.5142

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6091.295
$ws.Range("I32").Value = 2662.614
$ws.Range("J32").Value = 54950
$ws.Range("K32").Value = 2662.614
$ws.Range("L32").Value = 54950
$ws.Range("M32").Value = -2375.614
$ws.Range("N32").Value = -55524

# Row 102
$ws.Range("H102").Value = 1470.0588
$ws.Range("I102").Value = 1463.6428
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1463.6428
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 158.3571999999999
$ws.Range("N102").Value = -4744

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3248311.8
$ws.Range("I105").Value = 3789063.8
$ws.Range("J105").Value = 3800
$ws.Range("K105").Value = 3789063.8
$ws.Range("L105").Value = 3800
$ws.Range("M105").Value = -3787316.8
$ws.Range("N105").Value = -7294

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 398.75
$ws.Range("I22").Value = 305.91666
$ws.Range("J22").Value = 538
$ws.Range("K22").Value = 305.91666
$ws.Range("L22").Value = 538
$ws.Range("M22").Value = 44.08334000000002
$ws.Range("N22").Value = -1238

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 128.1
$ws.Range("I2").Value = 11.9375
$ws.Range("J2").Value = 260.85715
$ws.Range("K2").Value = 71.625
$ws.Range("L2").Value = 1565.1429
$ws.Range("M2").Value = 41.375
$ws.Range("N2").Value = -1791.1429

# Row 5
$ws.Range("H5").Value = 571.44446
$ws.Range("I5").Value = 477.57144
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 1432.71432
$ws.Range("L5").Value = 2700
$ws.Range("M5").Value = -1320.71432
$ws.Range("N5").Value = -2924

# Row 117
$ws.Range("H117").Value = 1315.1428
$ws.Range("I117").Value = 593
$ws.Range("J117").Value = 1512.091
$ws.Range("K117").Value = 1779
$ws.Range("L117").Value = 4536.272999999999
$ws.Range("M117").Value = 1663
$ws.Range("N117").Value = -11420.273

# Row 129
$ws.Range("H129").Value = 1531.5834
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1531.5834
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 4594.7502
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -14594.7502

# Row 135
$ws.Range("H135").Value = 571.44446
$ws.Range("I135").Value = 477.57144
$ws.Range("J135").Value = 900
$ws.Range("K135").Value = 4298.14296
$ws.Range("L135").Value = 8100
$ws.Range("M135").Value = -1763.14296
$ws.Range("N135").Value = -13170

# Row 137
$ws.Range("H137").Value = 2515.5
$ws.Range("I137").Value = 1072.3846
$ws.Range("J137").Value = 4221
$ws.Range("K137").Value = 3217.1538
$ws.Range("L137").Value = 12663
$ws.Range("M137").Value = 1882.8462
$ws.Range("N137").Value = -22863

$ws = $wb.Worksheets.Item("GSM")
# Row 110
$ws.Range("H110").Value = 43944
$ws.Range("J110").Value = 43944
$ws.Range("L110").Value = 43944
$ws.Range("N110").Value = -52124

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 886.73334
$ws.Range("I22").Value = 125.25
$ws.Range("J22").Value = 1163.6364
$ws.Range("K22").Value = 125.25
$ws.Range("L22").Value = 1163.6364
$ws.Range("M22").Value = 169.75
$ws.Range("N22").Value = -1753.6364

# Row 27
$ws.Range("H27").Value = 886.73334
$ws.Range("I27").Value = 125.25
$ws.Range("J27").Value = 1163.6364
$ws.Range("K27").Value = 125.25
$ws.Range("L27").Value = 1163.6364
$ws.Range("M27").Value = -18.25
$ws.Range("N27").Value = -1377.6364

# Row 132
$ws.Range("H132").Value = 2205.641
$ws.Range("I132").Value = 1491.6207
$ws.Range("K132").Value = 4474.8621
$ws.Range("M132").Value = -1944.8621

# Row 136
$ws.Range("H136").Value = 36623.035
$ws.Range("I136").Value = 50337.668
$ws.Range("J136").Value = 4622.222
$ws.Range("K136").Value = 151013.004
$ws.Range("L136").Value = 13866.666
$ws.Range("M136").Value = -148463.004
$ws.Range("N136").Value = -18966.666
